$wb = $excel.ActiveWorkbook

# The two "full" sheets (展览 / 全部类型) hold identical event-listing tables.
# Net effect of the update (per the diff): the oldest event row (the one
# that used to be row 2, "thp01~风摄少微") is removed, every subsequent
# event's B:I (date/name/venue/time/attendees/price/link/cover) data moves
# up one row, a couple of numbers tick up by one (想去人数 278 / 816), and
# the now-unused last row (row 8) is cleared out, shrinking the sheet's
# used range from A1:I8 down to A1:I7. Column A (the 0-based display
# index) already read 0,1,2,3,4,5,6,7 and needs no edits for rows 2-7.

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # B column holds dates formatted as plain "yyyy-mm-dd" text (not real
    # Excel dates) in the source file. Force text storage so assigning a
    # string like "2024-08-03" isn't auto-parsed into a date serial.
    $ws.Range("B2:B7").NumberFormat = "@"

    $ws.Range("B2").Value = "2024-08-03"
    $ws.Range("C2").Value = "丽水·樱卡动漫游戏嘉年华"
    $ws.Range("D2").Value = "中东路848号(解放街交汇) 飞达国际大酒店"
    $ws.Range("E2").Value = "2024.08.03 10:00-08.03 17:00"
    $ws.Range("F2").Value = 278
    $ws.Range("G2").Value = 50
    $ws.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=87276"
    $ws.Range("I2").Value = "//i0.hdslb.com/bfs/openplatform/202406/bVp0Zg1B1718172430380.jpeg"

    $ws.Range("B3").Value = "2024-08-03"
    $ws.Range("C3").Value = "丽水·逆光ZERO动漫游戏展"
    $ws.Range("D3").Value = "丽阳街651号 丽水华侨君澜大饭店"
    $ws.Range("E3").Value = "2024.08.03 10:00-08.03 17:00"
    $ws.Range("F3").Value = 6
    $ws.Range("G3").Value = 40
    $ws.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=89687"
    $ws.Range("I3").Value = "//i1.hdslb.com/bfs/openplatform/202407/NMYT1LRl1721639164353.jpeg"

    $ws.Range("B4").Value = "2024-08-10"
    $ws.Range("C4").Value = "丽水·CCAC动漫七夕（回馈展）"
    $ws.Range("D4").Value = "中东路848号(解放街交汇) 飞达国际大酒店"
    $ws.Range("E4").Value = "2024.08.10 09:00-08.10 17:00"
    $ws.Range("F4").Value = 92
    $ws.Range("G4").Value = 29.9
    $ws.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=86567"
    $ws.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202405/tsOzbBRx1717015539538.png"

    $ws.Range("B5").Value = "2024-08-17"
    $ws.Range("C5").Value = "丽水·AEO纯白礼赞动漫嘉年华"
    $ws.Range("D5").Value = "城北街1001号 爱依·时尚婚宴中心"
    $ws.Range("E5").Value = "2024.08.17 09:00-08.17 18:00"
    $ws.Range("F5").Value = 816
    $ws.Range("G5").Value = 55
    $ws.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=86779"
    $ws.Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202406/MxJ3oNjt1717405405850.jpeg"

    $ws.Range("B6").Value = "2024-08-24"
    $ws.Range("C6").Value = "丽水·R动漫嘉年华"
    $ws.Range("D6").Value = "中东路848号(解放街交汇) 飞达国际大酒店"
    $ws.Range("E6").Value = "2024.08.24 09:30-08.24 17:00"
    $ws.Range("F6").Value = 6
    $ws.Range("G6").Value = 45
    $ws.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=89651"
    $ws.Range("I6").Value = "//i0.hdslb.com/bfs/openplatform/202407/7o5ALbAM1721383424201.jpeg"

    $ws.Range("B7").Value = "2024-09-16"
    $ws.Range("C7").Value = "丽水·LZ栗子动漫游戏嘉年华"
    $ws.Range("D7").Value = "城北街798号 莱茵体育生活馆"
    $ws.Range("E7").Value = "2024.09.16 09:30-09.16 17:00"
    $ws.Range("F7").Value = 418
    $ws.Range("G7").Value = 65
    $ws.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=87480"
    $ws.Range("I7").Value = "//i1.hdslb.com/bfs/openplatform/202406/bATqcZhH1719285865931.jpeg"

    # Row 8 (old last event, "LZ栗子动漫游戏嘉年华") no longer exists; clear
    # it entirely so the used range shrinks to A1:I7.
    $ws.Range("A8:I8").Clear()
}
